$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "600.55")
# but must be preserved as plain text, matching the source data which
# stores every Price cell as an inline string. Force Text format before
# assigning so Excel does not silently convert these into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.510.32"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.511.57"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.55"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.98"
$ws.Range("E6").Value = "  +4.41%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.132"
$ws.Range("E9").Value = "  -0.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.432"
$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.111.14"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "31.33"
$ws.Range("E13").Value = "  +10.86%  "

$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.458.50"
$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000180"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.491.28"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.32"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.66"
$ws.Range("E19").Value = "  +4.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.25"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.05"
$ws.Range("E21").Value = "  +1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.55"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.540"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.72"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.28"
$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.993"
$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.16"
$ws.Range("E30").Value = "  -2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -1.79%  "

$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.79"
$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.41"
$ws.Range("E34").Value = "  +0.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.65"
$ws.Range("E35").Value = "  +2.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.26"
$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  +2.46%  "

$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.06"
$ws.Range("E39").Value = "  +3.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.70"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.48"
$ws.Range("E41").Value = "  +2.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0736"
$ws.Range("E42").Value = "  -0.86%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.31"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.814.97"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.60"
$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0302"
$ws.Range("E47").Value = "  -2.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "340.15"
$ws.Range("E48").Value = "  -0.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.13"
$ws.Range("E50").Value = "  +1.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.852"
$ws.Range("E51").Value = "  -0.19%  "
